$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns to snake_case field names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Title-case the connector words (de/del/la/las/el/los/y) in state/municipality names
$ws.Range("B6").Value = "Rincón De Romos"
$ws.Range("B7").Value = "San José De Gracia"
$ws.Range("B22").Value = "Amatenango De La Frontera"
$ws.Range("B24").Value = "Bejucal De Ocampo"
$ws.Range("B30").Value = "Comitán De Domínguez"
$ws.Range("B45").Value = "Mazapa De Madero"
$ws.Range("B47").Value = "Montecristo De Guerero"
$ws.Range("B54").Value = "San Cristóbal De Las Casas"
$ws.Range("B86").Value = "Villa De Álvarez"
$ws.Range("A88").Value = "Ciudad De México"
$ws.Range("B102").Value = "Coneto De Comonfort"
$ws.Range("A116").Value = "Estado De México"
$ws.Range("B116").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B117").Value = "Almoloya De Alquisiras"
$ws.Range("B118").Value = "Almoloya De Juárez"
$ws.Range("B119").Value = "Almoloya Del Río"
$ws.Range("B128").Value = "Chapa De Mota"
$ws.Range("B133").Value = "Ecatepec De Morelos"
$ws.Range("B135").Value = "Ixtapan De La Sal"
$ws.Range("B144").Value = "Naucalpan De Juárez"
$ws.Range("B149").Value = "San Felipe Del Progreso"
$ws.Range("B150").Value = "Soyaniquilpan De Juárez"
$ws.Range("B158").Value = "Tenango Del Valle"
$ws.Range("B161").Value = "Tlalnepantla De Baz"
$ws.Range("B165").Value = "Valle De Bravo"
$ws.Range("A168").Value = "Guanajuato"
$ws.Range("B176").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B183").Value = "Purísima Del Rincón"
$ws.Range("B186").Value = "San Francisco Del Rincón"
$ws.Range("B188").Value = "San Luis De La Paz"
$ws.Range("B189").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B190").Value = "Silao De La Victoria"
$ws.Range("B197").Value = "Acapulco De Juárez"
$ws.Range("B200").Value = "Ayutla De Los Libres"
$ws.Range("B202").Value = "Buenavista De Cuéllar"
$ws.Range("B203").Value = "Chilapa De Álvarez"
$ws.Range("B204").Value = "Chilpancingo De Los Bravo"
$ws.Range("B206").Value = "Coyuca De Benítez"
$ws.Range("B207").Value = "Coyuca De Catalán"
$ws.Range("B209").Value = "Cutzamala De Pinzón"
$ws.Range("B212").Value = "Huitzuco De Los Figueroa"
$ws.Range("B213").Value = "Iguala De La Independencia"
$ws.Range("B215").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B216").Value = "Mártir De Cuilapan"
$ws.Range("B223").Value = "Taxco De Alarcón"
$ws.Range("B224").Value = "Técpan De Galeana"
$ws.Range("B226").Value = "Tixtla De Guerero"
$ws.Range("B230").Value = "Tlapa De Comonfort"
$ws.Range("B235").Value = "Atotonilco El Grande"
$ws.Range("B238").Value = "Cuautepec De Hinojosa"
$ws.Range("B240").Value = "Huasca De Ocampo"
$ws.Range("B242").Value = "Huejutla De Reyes"
$ws.Range("B246").Value = "Mixquiahuala De Juárez"
$ws.Range("B248").Value = "Pachuca De Soto"
$ws.Range("B253").Value = "Tula De Allende"
$ws.Range("B254").Value = "Tulancingo De Bravo"
$ws.Range("B260").Value = "Atemajac De Brizuela"
$ws.Range("B261").Value = "Atotonilco El Alto"
$ws.Range("B262").Value = "Autlán De Navarro"
$ws.Range("B268").Value = "Encarnación De Díaz"
$ws.Range("B273").Value = "Jilotlán De Los Dolores"
$ws.Range("B275").Value = "Lagos De Moreno"
$ws.Range("B277").Value = "San Diego De Alejandría"
$ws.Range("B285").Value = "Unión De San Antonio"
$ws.Range("B329").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B350").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B352").Value = "Puente De Ixtla"
$ws.Range("B355").Value = "Tetela Del Volcán"
$ws.Range("B356").Value = "Tlaltizapán De Zapata"
$ws.Range("B359").Value = "Zacualpan De Amilpas"
$ws.Range("B372").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B373").Value = "Ayoquezco De Aldama"
$ws.Range("B377").Value = "Coicoyán De Las Flores"
$ws.Range("B380").Value = "El Barrio De La Soledad"
$ws.Range("B381").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B382").Value = "Ixtlán De Juárez"
$ws.Range("B383").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B388").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B389").Value = "Oaxaca De Juárez"
$ws.Range("B390").Value = "Ocotlán De Morelos"
$ws.Range("B391").Value = "Putla Villa De Guerero"
$ws.Range("B397").Value = "San Antonino El Alto"
$ws.Range("B441").Value = "Santa María Jalapa Del Marqués"
$ws.Range("B460").Value = "Santo Domingo De Morelos"
$ws.Range("B465").Value = "Tamazulápam Del Espíritu Santo"
$ws.Range("B466").Value = "Tezoatlán De Segura Y Luna"
$ws.Range("B467").Value = "Tlacolula De Matamoros"
$ws.Range("B468").Value = "Totontepec Villa De Morelos"
$ws.Range("B469").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B470").Value = "Zimatlán De Álvarez"
$ws.Range("B482").Value = "Ixcamilpa De Guerero"
$ws.Range("B483").Value = "Izúcar De Matamoros"
$ws.Range("B487").Value = "Palmar De Bravo"
$ws.Range("B499").Value = "Tepexi De Rodríguez"
$ws.Range("B502").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B513").Value = "Amealco De Bonfil"
$ws.Range("B515").Value = "Cadereyta De Montes"
$ws.Range("B517").Value = "Landa De Matamoros"
$ws.Range("B518").Value = "Pinal De Amoles"
$ws.Range("B526").Value = "Axtla De Terrazas"
$ws.Range("B533").Value = "Mexquitic De Carmona"
$ws.Range("B538").Value = "Santa María Del Río"
$ws.Range("B541").Value = "Villa De Arista"
$ws.Range("B542").Value = "Villa De Guadalupe"
$ws.Range("B543").Value = "Villa De Ramos"
$ws.Range("B557").Value = "San Miguel De Horcasitas"
$ws.Range("B577").Value = "Soto La Marina"
$ws.Range("B592").Value = "Alto Lucero De Gutiérrez Barrios"
$ws.Range("B594").Value = "Amatlán De Los Reyes"
$ws.Range("B597").Value = "Boca Del Río"
$ws.Range("B606").Value = "Cosamaloapan De Carpio"
$ws.Range("B613").Value = "Hueyapan De Ocampo"
$ws.Range("B615").Value = "Ixhuatlán De Madero"
$ws.Range("B620").Value = "Juchique De Ferrer"
$ws.Range("B625").Value = "Martínez De La Torre"
$ws.Range("B629").Value = "Nanchital De Lázaro Cárdenas Del Río"
$ws.Range("B637").Value = "Poza Rica De Hidalgo"
$ws.Range("B643").Value = "Sayula De Alemán"
$ws.Range("B644").Value = "Soledad De Doblado"
$ws.Range("B647").Value = "Tatahuicapan De Juárez"
$ws.Range("B655").Value = "Tlacotepec De Mejía"
$ws.Range("B659").Value = "Vega De Alatorre"
$ws.Range("B672").Value = "Noria De Ángeles"

# Remove trailing footnote/source rows (682-686); this also shrinks the used range / dimension to A1:D680
$ws.Range("A682:D686").EntireRow.Delete()
